# Fruta / hortaliza, semanal
#
# Inserts two new weekly price-report rows (date 44461 = 2021-09-22) for
# "Femacal de La Calera" / Coquimbo / Brocoli, one per quality grade
# ("Primera" and "Segunda"), at the top of that market's data block
# (rows 350-351), pushing the existing rows 350-377 down to 352-379.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing block (old rows 350:377) down by two rows so the
# two brand-new rows can be inserted at the top of the block.
$ws.Rows("350:351").Insert()

# New row 350: Primera
$ws.Range("A350").Value = 3
$ws.Range("B350").Value = "Femacal de La Calera"
$ws.Range("C350").Value = "Coquimbo"
$ws.Range("D350").Value = 44461
$ws.Range("E350").Value = 5
$ws.Range("F350").Value = 100112023
$ws.Range("G350").Value = "Brócoli"
$ws.Range("H350").Value = "Sin especificar"
$ws.Range("I350").Value = "Primera"
$ws.Range("J350").Value = 2250
$ws.Range("K350").Value = 550
$ws.Range("L350").Value = 600
$ws.Range("M350").Value = 571
$ws.Range("N350").Value = "`$/unidad"
$ws.Range("O350").Value = "Provincia de Quillota"
$ws.Range("P350").Value = 571
$ws.Range("Q350").Value = 1
$ws.Range("R350").Value = "Hortaliza"

# New row 351: Segunda
$ws.Range("A351").Value = 3
$ws.Range("B351").Value = "Femacal de La Calera"
$ws.Range("C351").Value = "Coquimbo"
$ws.Range("D351").Value = 44461
$ws.Range("E351").Value = 5
$ws.Range("F351").Value = 100112023
$ws.Range("G351").Value = "Brócoli"
$ws.Range("H351").Value = "Sin especificar"
$ws.Range("I351").Value = "Segunda"
$ws.Range("J351").Value = 900
$ws.Range("K351").Value = 450
$ws.Range("L351").Value = 450
$ws.Range("M351").Value = 450
$ws.Range("N351").Value = "`$/unidad"
$ws.Range("O351").Value = "Provincia de Quillota"
$ws.Range("P351").Value = 450
$ws.Range("Q351").Value = 1
$ws.Range("R351").Value = "Hortaliza"
